$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$fmt = $ws.Range("F3").NumberFormat
Write-Output "fmt=$fmt"
$ws.Range("Z1").Value2 = [datetime]"2021-04-30"
$ws.Range("Z1").NumberFormat = $fmt
